# Updated cryptos list (price + 1h volume change) - generated from source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.978.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.87%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4597"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07733"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9870"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.944.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.686"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07050"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "84.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009541"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.977.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.334"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.152.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.610"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.837"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09268"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.257"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.016"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05724"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02045"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5538"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.298"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.748"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.084"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06825"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.780"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.44%  "
